$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404"
$headers_FV2310 = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

$headers_FV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers_FV2310.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers_FV2310[$i]
}

# Column K (11) is "diff" - stays the same

for ($i = 0; $i -lt $headers_FV2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headers_FV2404[$i]
}

# Turn the used range into an Excel table ("Table1")
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U82"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1)
$ws.Activate()
$ws.Cells.Item(2, 1).Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$null
